$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E2").Value = 11

$ws.Range("D3").Value = 10.71

$ws.Range("C4").Value = 9.289999999999999
$ws.Range("E4").Value = 9.92
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 10

$ws.Range("B5").Value = 9
$ws.Range("D5").Value = 10.08
$ws.Range("F5").Value = 10.27
$ws.Range("G5").Value = 9.32

$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 9.73
$ws.Range("H6").Value = 11.42

$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 10.68
$ws.Range("H7").Value = 9.83

$ws.Range("F8").Value = 8.58
$ws.Range("G8").Value = 10.17
